$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Towards a responsible innovation agenda for HCI", "['Samantha']"),
    @("SpaceBot: Towards participatory evaluation of smart buildings", "['Samantha']"),
    @("Augmenting Audits", "['Samantha']"),
    @("Analysis of lysosomal enzyme activities in  induced pluripotent stem cell, neural progenitor  cell, and neuron models as potential biomarkers  of Huntington’s Disease", "['Callum']"),
    @("Assessing the role of complement domain containing peptides in the developing fly brain", "['Callum']"),
    @("Analysis of lysosomal enzyme activities in  induced pluripotent stem cell, neural progenitor  cell, and neuron models as potential biomarkers  of Huntington’s Disease", "['Markus']")
)

$row = 6
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
